# Auto-generated: update cached market-data values per upstream scheduled-runner refresh.
# Source: diff against Sheets/Tonberry_Profits.xlsx (concatenated sheet export).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2998.75
$ws.Range("I40").Value = 2998.5
$ws.Range("J40").Value = 2999
$ws.Range("K40").Value = 2998.5
$ws.Range("L40").Value = 2999
$ws.Range("M40").Value = -2823.5
$ws.Range("N40").Value = -3349
$ws.Range("H43").Value = 1582
$ws.Range("J43").Value = 1709.5
$ws.Range("L43").Value = 1709.5
$ws.Range("N43").Value = -1847.5
$ws.Range("H95").Value = 28666
$ws.Range("J95").Value = 28666
$ws.Range("L95").Value = 28666
$ws.Range("N95").Value = -34158
$ws.Range("H101").Value = 650
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -6244
$ws.Range("H111").Value = 645
$ws.Range("I111").Value = 645
$ws.Range("K111").Value = 1935
$ws.Range("M111").Value = 1132
$ws.Range("H113").Value = 27274.625
$ws.Range("I113").Value = 35416.668
$ws.Range("J113").Value = 2848.5
$ws.Range("K113").Value = 35416.668
$ws.Range("L113").Value = 2848.5
$ws.Range("M113").Value = -32162.668
$ws.Range("N113").Value = -9356.5
$ws.Range("H131").Value = 2053.85
$ws.Range("J131").Value = 3410.7
$ws.Range("L131").Value = 10232.1
$ws.Range("N131").Value = -20312.1
$ws.Range("H140").Value = 47929.062
$ws.Range("J140").Value = 47929.062
$ws.Range("L140").Value = 47929.062
$ws.Range("N140").Value = -58289.062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2879.362
$ws.Range("I32").Value = 1947.9574
$ws.Range("K32").Value = 1947.9574
$ws.Range("M32").Value = -1660.9574
$ws.Range("H45").Value = 1444.6
$ws.Range("I45").Value = 884.1111
$ws.Range("J45").Value = 2885.8572
$ws.Range("K45").Value = 884.1111
$ws.Range("L45").Value = 2885.8572
$ws.Range("M45").Value = -507.1111
$ws.Range("N45").Value = -3639.8572
$ws.Range("H74").Value = 4900
$ws.Range("I74").Value = 4900
$ws.Range("K74").Value = 4900
$ws.Range("M74").Value = -4026
$ws.Range("H77").Value = 4900
$ws.Range("I77").Value = 4900
$ws.Range("K77").Value = 24500
$ws.Range("M77").Value = -20132
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H102").Value = 874.75
$ws.Range("I102").Value = 874.75
$ws.Range("K102").Value = 874.75
$ws.Range("M102").Value = 747.25
$ws.Range("H109").Value = 67753.28999999999
$ws.Range("J109").Value = 67753.28999999999
$ws.Range("L109").Value = 67753.28999999999
$ws.Range("N109").Value = -70527.28999999999
$ws.Range("H110").Value = 1760.3572
$ws.Range("I110").Value = 1449.5454
$ws.Range("J110").Value = 2900
$ws.Range("K110").Value = 1449.5454
$ws.Range("L110").Value = 2900
$ws.Range("M110").Value = 595.4546
$ws.Range("N110").Value = -6990
$ws.Range("H122").Value = 1571.2
$ws.Range("I122").Value = 1541.0834
$ws.Range("K122").Value = 4623.2502
$ws.Range("M122").Value = -2173.2502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 402400
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 668666.7
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 668666.7
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -670912.7
$ws.Range("H89").Value = 402400
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 668666.7
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 3343333.5
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -3354565.5
$ws.Range("H94").Value = 936.0714
$ws.Range("I94").Value = 766.2222
$ws.Range("K94").Value = 766.2222
$ws.Range("M94").Value = -315.2222
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H107").Value = 2904.2354
$ws.Range("J107").Value = 3741.6667
$ws.Range("L107").Value = 3741.6667
$ws.Range("N107").Value = -7581.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1660
$ws.Range("J22").Value = 1660
$ws.Range("L22").Value = 1660
$ws.Range("N22").Value = -2360
$ws.Range("H31").Value = 3012.8572
$ws.Range("I31").Value = 1239.4375
$ws.Range("K31").Value = 1239.4375
$ws.Range("M31").Value = -944.4375
$ws.Range("H34").Value = 3012.8572
$ws.Range("I34").Value = 1239.4375
$ws.Range("K34").Value = 1239.4375
$ws.Range("M34").Value = -1037.4375
$ws.Range("H86").Value = 1650
$ws.Range("I86").Value = 1650
$ws.Range("K86").Value = 1650
$ws.Range("M86").Value = -527
$ws.Range("H89").Value = 1650
$ws.Range("I89").Value = 1650
$ws.Range("K89").Value = 8250
$ws.Range("M89").Value = -2634
$ws.Range("H92").Value = 30500
$ws.Range("J92").Value = 30500
$ws.Range("L92").Value = 30500
$ws.Range("N92").Value = -35492
$ws.Range("H106").Value = 42900
$ws.Range("J106").Value = 42900
$ws.Range("L106").Value = 42900
$ws.Range("N106").Value = -45424

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 41191.445
$ws.Range("I129").Value = 771.1667
$ws.Range("J129").Value = 61401.582
$ws.Range("K129").Value = 2313.5001
$ws.Range("L129").Value = 184204.746
$ws.Range("M129").Value = 2686.4999
$ws.Range("N129").Value = -194204.746
$ws.Range("H131").Value = 13911757
$ws.Range("I131").Value = 83333840
$ws.Range("J131").Value = 27340.633
$ws.Range("K131").Value = 250001520
$ws.Range("L131").Value = 82021.899
$ws.Range("M131").Value = -249996480
$ws.Range("N131").Value = -92101.899

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984
$ws.Range("H97").Value = 1311.2307
$ws.Range("I97").Value = 444.83334
$ws.Range("J97").Value = 2053.8572
$ws.Range("K97").Value = 444.83334
$ws.Range("L97").Value = 2053.8572
$ws.Range("M97").Value = 51.16665999999998
$ws.Range("N97").Value = -3045.8572
$ws.Range("H101").Value = 14250
$ws.Range("J101").Value = 14250
$ws.Range("L101").Value = 14250
$ws.Range("N101").Value = -20740
$ws.Range("H102").Value = 1851.2727
$ws.Range("I102").Value = 1739.4546
$ws.Range("J102").Value = 1963.091
$ws.Range("K102").Value = 1739.4546
$ws.Range("L102").Value = 1963.091
$ws.Range("M102").Value = -117.4546
$ws.Range("N102").Value = -5207.091
$ws.Range("H122").Value = 1798.5714
$ws.Range("I122").Value = 1644.8422
$ws.Range("K122").Value = 4934.5266
$ws.Range("M122").Value = -2484.5266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1360.1428
$ws.Range("I22").Value = 1110.75
$ws.Range("J22").Value = 1459.9
$ws.Range("K22").Value = 1110.75
$ws.Range("L22").Value = 1459.9
$ws.Range("M22").Value = -815.75
$ws.Range("N22").Value = -2049.9
$ws.Range("H27").Value = 1360.1428
$ws.Range("I27").Value = 1110.75
$ws.Range("J27").Value = 1459.9
$ws.Range("K27").Value = 1110.75
$ws.Range("L27").Value = 1459.9
$ws.Range("M27").Value = -1003.75
$ws.Range("N27").Value = -1673.9
$ws.Range("H46").Value = 1791.5834
$ws.Range("J46").Value = 1791.5834
$ws.Range("L46").Value = 1791.5834
$ws.Range("N46").Value = -2167.5834
$ws.Range("H136").Value = 4919.933
$ws.Range("J136").Value = 5999.875
$ws.Range("L136").Value = 17999.625
$ws.Range("N136").Value = -23099.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 494.5
$ws.Range("I81").Value = 494.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 989
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 72
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 494.5
$ws.Range("I84").Value = 494.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4945
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 359
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 800.1429000000001
$ws.Range("H122").Value = 21010.924
$ws.Range("I122").Value = 29421.777
$ws.Range("K122").Value = 88265.33099999999
$ws.Range("M122").Value = -85815.33099999999
